$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1497.8889
$ws.Range("I18").Value = 1497.8889
$ws.Range("K18").Value = 1497.8889
$ws.Range("M18").Value = -1213.8889
$ws.Range("H43").Value = 11970.667
$ws.Range("I43").Value = 17362.5
$ws.Range("J43").Value = 1187
$ws.Range("K43").Value = 17362.5
$ws.Range("L43").Value = 1187
$ws.Range("M43").Value = -17293.5
$ws.Range("N43").Value = -1325
$ws.Range("H62").Value = 6859.4443
$ws.Range("I62").Value = 7390.857
$ws.Range("J62").Value = 4999.5
$ws.Range("K62").Value = 7390.857
$ws.Range("L62").Value = 4999.5
$ws.Range("M62").Value = -6766.857
$ws.Range("N62").Value = -6247.5
$ws.Range("H65").Value = 6859.4443
$ws.Range("I65").Value = 7390.857
$ws.Range("J65").Value = 4999.5
$ws.Range("K65").Value = 36954.285
$ws.Range("L65").Value = 24997.5
$ws.Range("M65").Value = -33834.285
$ws.Range("N65").Value = -31237.5
$ws.Range("H87").Value = 70000
$ws.Range("I87").Value = 70000
$ws.Range("K87").Value = 70000
$ws.Range("M87").Value = -68752
$ws.Range("H90").Value = 70000
$ws.Range("I90").Value = 70000
$ws.Range("K90").Value = 210000
$ws.Range("M90").Value = -203760
$ws.Range("H106").Value = 1665.6666
$ws.Range("I106").Value = 1665.6666
$ws.Range("K106").Value = 1665.6666
$ws.Range("M106").Value = -1034.6666
$ws.Range("H133").Value = 78987.5
$ws.Range("J133").Value = 78987.5
$ws.Range("L133").Value = 78987.5
$ws.Range("N133").Value = -89107.5
$ws.Range("H137").Value = 7140.722
$ws.Range("I137").Value = 1628.9333
$ws.Range("J137").Value = 34699.668
$ws.Range("K137").Value = 4886.7999
$ws.Range("L137").Value = 104099.004
$ws.Range("M137").Value = -2336.7999
$ws.Range("N137").Value = -109199.004
$ws.Range("H138").Value = 22167.89
$ws.Range("I138").Value = 51851.5
$ws.Range("J138").Value = 5205.8286
$ws.Range("K138").Value = 155554.5
$ws.Range("L138").Value = 15617.4858
$ws.Range("M138").Value = -150414.5
$ws.Range("N138").Value = -25897.4858

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 150061.75
$ws.Range("I45").Value = 196665.67
$ws.Range("J45").Value = 10250
$ws.Range("K45").Value = 196665.67
$ws.Range("L45").Value = 10250
$ws.Range("M45").Value = -196288.67
$ws.Range("N45").Value = -11004
$ws.Range("H74").Value = 33430.758
$ws.Range("I74").Value = 39819.117
$ws.Range("K74").Value = 39819.117
$ws.Range("M74").Value = -38945.117
$ws.Range("H77").Value = 33430.758
$ws.Range("I77").Value = 39819.117
$ws.Range("K77").Value = 199095.585
$ws.Range("M77").Value = -194727.585
$ws.Range("H122").Value = 1847.1428
$ws.Range("I122").Value = 1653.7222
$ws.Range("K122").Value = 4961.1666
$ws.Range("M122").Value = -2511.1666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2143.2222
$ws.Range("I105").Value = 1755.5714
$ws.Range("J105").Value = 3500
$ws.Range("K105").Value = 1755.5714
$ws.Range("L105").Value = 3500
$ws.Range("M105").Value = -8.57140000000004
$ws.Range("N105").Value = -6994
$ws.Range("H132").Value = 101814.836
$ws.Range("J132").Value = 101814.836
$ws.Range("L132").Value = 101814.836
$ws.Range("N132").Value = -111934.836
$ws.Range("H135").Value = 97484.5
$ws.Range("J135").Value = 97484.5
$ws.Range("L135").Value = 97484.5
$ws.Range("N135").Value = -107624.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 36326.2
$ws.Range("I31").Value = 54669.95
$ws.Range("K31").Value = 54669.95
$ws.Range("M31").Value = -54374.95
$ws.Range("H34").Value = 36326.2
$ws.Range("I34").Value = 54669.95
$ws.Range("K34").Value = 54669.95
$ws.Range("M34").Value = -54467.95
$ws.Range("H99").Value = 4999.625
$ws.Range("I99").Value = 2816.8
$ws.Range("J99").Value = 8637.666999999999
$ws.Range("K99").Value = 2816.8
$ws.Range("L99").Value = 8637.666999999999
$ws.Range("M99").Value = -1318.8
$ws.Range("N99").Value = -11633.667
$ws.Range("H105").Value = 952
$ws.Range("I105").Value = 868
$ws.Range("J105").Value = 1750
$ws.Range("K105").Value = 868
$ws.Range("L105").Value = 1750
$ws.Range("M105").Value = 879
$ws.Range("N105").Value = -5244
$ws.Range("H120").Value = 47714.145
$ws.Range("J120").Value = 47714.145
$ws.Range("L120").Value = 47714.145
$ws.Range("N120").Value = -54972.145
$ws.Range("H126").Value = 4999.625
$ws.Range("I126").Value = 2816.8
$ws.Range("J126").Value = 8637.666999999999
$ws.Range("K126").Value = 8450.400000000001
$ws.Range("L126").Value = 25913.001
$ws.Range("M126").Value = -5980.400000000001
$ws.Range("N126").Value = -30853.001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 89.72727
$ws.Range("I10").Value = 117.125
$ws.Range("J10").Value = 16.666666
$ws.Range("K10").Value = 351.375
$ws.Range("L10").Value = 49.999998
$ws.Range("M10").Value = -212.375
$ws.Range("N10").Value = -327.999998
$ws.Range("H60").Value = 1344.6428
$ws.Range("I60").Value = 1534.8334
$ws.Range("J60").Value = 203.5
$ws.Range("K60").Value = 4604.5002
$ws.Range("L60").Value = 610.5
$ws.Range("M60").Value = -4353.5002
$ws.Range("N60").Value = -1112.5
$ws.Range("H76").Value = 1132.5
$ws.Range("I76").Value = 1132.5
$ws.Range("K76").Value = 3397.5
$ws.Range("M76").Value = -3014.5
$ws.Range("H79").Value = 1132.5
$ws.Range("I79").Value = 1132.5
$ws.Range("K79").Value = 3397.5
$ws.Range("M79").Value = -2071.5
$ws.Range("H138").Value = 50010548
$ws.Range("I138").Value = 55565612
$ws.Range("J138").Value = 15000
$ws.Range("K138").Value = 166696836
$ws.Range("L138").Value = 45000
$ws.Range("M138").Value = -166691696
$ws.Range("N138").Value = -55280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8262
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("H83").Value = 8262
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("H131").Value = 99497
$ws.Range("J131").Value = 99497
$ws.Range("L131").Value = 99497
$ws.Range("N131").Value = -109577
$ws.Range("H132").Value = 3873.1143
$ws.Range("J132").Value = 4328.5
$ws.Range("L132").Value = 12985.5
$ws.Range("N132").Value = -18045.5
$ws.Range("M80").ClearContents()
$ws.Range("M83").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1118.625
$ws.Range("I46").Value = 1091.5
$ws.Range("J46").Value = 1200
$ws.Range("K46").Value = 1091.5
$ws.Range("L46").Value = 1200
$ws.Range("M46").Value = -903.5
$ws.Range("N46").Value = -1576
$ws.Range("H55").Value = 617.6
$ws.Range("I55").Value = 415.27274
$ws.Range("K55").Value = 415.27274
$ws.Range("M55").Value = -242.27274
$ws.Range("H76").Value = 10048.667
$ws.Range("J76").Value = 11507.75
$ws.Range("L76").Value = 11507.75
$ws.Range("N76").Value = -12183.75
$ws.Range("H79").Value = 10048.667
$ws.Range("J79").Value = 11507.75
$ws.Range("L79").Value = 11507.75
$ws.Range("N79").Value = -13847.75
$ws.Range("H122").Value = 561026.1
$ws.Range("I122").Value = 838064.5600000001
$ws.Range("K122").Value = 2514193.68
$ws.Range("M122").Value = -2511743.68
$ws.Range("H138").Value = 54384.5
$ws.Range("J138").Value = 78769
$ws.Range("L138").Value = 78769
$ws.Range("N138").Value = -89049

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6299.3335
$ws.Range("J62").Value = 6003
$ws.Range("L62").Value = 6003
$ws.Range("N62").Value = -7251
$ws.Range("H65").Value = 6299.3335
$ws.Range("J65").Value = 6003
$ws.Range("L65").Value = 30015
$ws.Range("N65").Value = -36255
$ws.Range("H132").Value = 1219.2413
$ws.Range("I132").Value = 1235.6666
$ws.Range("K132").Value = 3706.9998
$ws.Range("M132").Value = -1176.9998
$ws.Range("H136").Value = 2611.2258
$ws.Range("I136").Value = 2423.3
$ws.Range("J136").Value = 2952.9092
$ws.Range("K136").Value = 7269.900000000001
$ws.Range("L136").Value = 8858.7276
$ws.Range("M136").Value = -4719.900000000001
$ws.Range("N136").Value = -13958.7276
$ws.Range("H140").Value = 77197.25
$ws.Range("J140").Value = 77197.25
$ws.Range("L140").Value = 77197.25
$ws.Range("N140").Value = -87557.25
